# TimeTracking.xlsx — fill in the logged hours for Week 7 (row 9) and
# Week 8 (row 10), which lets the running "Total" columns (C, E, G, I)
# pick up real values instead of carrying the prior week's total forward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Week 7 (row 9): Hours entered for each of the four members ---
$ws.Range("B9").Value = 13
$ws.Range("D9").Value = 11
$ws.Range("F9").Value = 12
$ws.Range("H9").Value = 11

# C9/E9 already carry the running-total shared formula down from C8/E8,
# so they recompute automatically. G9/I9 only had the trivial "=F9"/"=H9"
# placeholder (because the row used to be blank) - now that there's real
# data, give them the same running-total pattern the rows above use.
$ws.Range("G9").Formula = "=F9+G8"
$ws.Range("I9").Formula = "=H9+I8"

# --- Week 8 (row 10): Hours entered for each of the four members ---
$ws.Range("B10").Value = 12
$ws.Range("D10").Value = 14
$ws.Range("F10").Value = 11
$ws.Range("H10").Value = 13

$ws.Range("G10").Formula = "=F10+G9"
$ws.Range("I10").Formula = "=H10+I9"

# --- Leave the workbook with H10 selected, matching where editing stopped ---
$ws.Range("H10").Select() | Out-Null
